# Updated test data for German,Czech market
# Added test data for Belgium market
#
# Removes the "ZXF" and "ZXFEV" rows (rows 20 and 21) from the Germany and
# Belgium worksheets, since those repeater models no longer apply. This
# shifts the "Wg" / "Repeaters" footer rows up to rows 20/21.

$wb = $excel.ActiveWorkbook

$germany = $wb.Worksheets.Item("Germany")
$belgium = $wb.Worksheets.Item("Belgium")

# Delete rows 20 and 21 (ZXF, ZXFEV) on the Germany sheet; rows below
# shift up so "Wg"/"Repeaters" become rows 20/21.
[void]$germany.Range("A20:D21").EntireRow.Delete()

# Delete rows 20 and 21 (ZXF, ZXFEV) on the Belgium sheet.
[void]$belgium.Range("A20:D21").EntireRow.Delete()

# Restore the selections recorded in the saved workbook.
[void]$germany.Activate()
[void]$germany.Range("A14").Select()

[void]$belgium.Activate()
[void]$belgium.Range("A20:XFD21").Select()

# The Czech sheet (tab index 2) remains the active tab, as in the saved file.
[void]$wb.Worksheets.Item("Czech").Activate()
